# Fruta / hortaliza, semanal
#
# A new weekly price record was inserted for
# "Feria Lagunitas de Puerto Montt - Frutilla" at row 494, pushing every
# existing record from row 494 down one row (494->495, ..., 523->524) and
# growing the used range from A1:T523 to A1:T524.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 494; this shifts rows
# 494..523 down to 495..524 and widens the sheet dimension automatically.
$ws.Rows.Item(494).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A494").Value = 4
$ws.Range("B494").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C494").Value = "Los Lagos"
$ws.Range("D494").Value = 45267
$ws.Range("E494").Value = 10
$ws.Range("F494").Value = "Fruta"
$ws.Range("G494").Value = 100101
$ws.Range("H494").Value = "Berries"
$ws.Range("I494").Value = 100112025
$ws.Range("J494").Value = "Frutilla"
$ws.Range("K494").Value = "Sin especificar"
$ws.Range("L494").Value = "Primera"
$ws.Range("M494").Value = 300
$ws.Range("N494").Value = 15000
$ws.Range("O494").Value = 15000
$ws.Range("P494").Value = 15000
$ws.Range("Q494").Value = "$/bandeja 7 kilos"
$ws.Range("R494").Value = "Provincia de Melipilla"
$ws.Range("S494").Value = 2143
$ws.Range("T494").Value = 7
